$d = $word.ActiveDocument

$replacements = @(
    @("79÷6=", "57÷4="),
    @("14÷6=", "88÷8="),
    @("57÷7=", "80÷6="),
    @("26÷3=", "48÷9="),
    @("76÷7=", "90÷4="),
    @("60÷8=", "83÷4="),
    @("27÷4=", "14÷5="),
    @("66÷3=", "20÷3="),
    @("20÷2=", "68÷3="),
    @("78÷3=", "82÷6="),
    @("89÷9=", "97÷5="),
    @("29÷4=", "11÷6="),
    @("18÷2=", "32÷4="),
    @("12÷3=", "28÷3="),
    @("12÷4=", "42÷2="),
    @("19÷3=", "49÷9="),
    @("64÷4=", "79÷2="),
    @("97÷6=", "14÷7="),
    @("54÷3=", "37÷9="),
    @("58÷7=", "38÷8="),
    @("33÷6=", "95÷5="),
    @("66÷4=", "19÷9="),
    @("10÷8=", "76÷4="),
    @("74÷4=", "13÷6="),
    @("73÷6=", "71÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
